$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Vorname" (first name) column B and "Nachname" (last name) column C
# had been mapped to the wrong shared-string values; swap the two columns
# (including the header row) so each row's values line up with the correct
# header again.
$rowCount = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $c
    $ws.Cells.Item($r, 3).Value = $b
}

$ws.Range("E4").Select()
